$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 4).Value = '29.083.14'
$ws.Cells.Item(2, 5).Value = '  -0.47%  '

$ws.Cells.Item(3, 4).Value = '1.821.02'
$ws.Cells.Item(3, 5).Value = '  -0.45%  '

$ws.Cells.Item(4, 5).Value = '  +0.45%  '

$ws.Cells.Item(5, 4).NumberFormat = '@'
$ws.Cells.Item(5, 4).Value = '232.98'
$ws.Cells.Item(5, 5).Value = '  -2.02%  '

$ws.Cells.Item(6, 4).NumberFormat = '@'
$ws.Cells.Item(6, 4).Value = '0.5901'
$ws.Cells.Item(6, 5).Value = '  -3.17%  '

$ws.Cells.Item(7, 4).NumberFormat = '@'
$ws.Cells.Item(7, 4).Value = '1.005'
$ws.Cells.Item(7, 5).Value = '  +0.41%  '

$ws.Cells.Item(8, 4).NumberFormat = '@'
$ws.Cells.Item(8, 4).Value = '0.2742'
$ws.Cells.Item(8, 5).Value = '  -3.07%  '

$ws.Cells.Item(9, 4).NumberFormat = '@'
$ws.Cells.Item(9, 4).Value = '0.06790'
$ws.Cells.Item(9, 5).Value = '  -4.44%  '

$ws.Cells.Item(10, 4).NumberFormat = '@'
$ws.Cells.Item(10, 4).Value = '23.00'
$ws.Cells.Item(10, 5).Value = '  -4.09%  '

$ws.Cells.Item(11, 4).NumberFormat = '@'
$ws.Cells.Item(11, 4).Value = '0.07513'
$ws.Cells.Item(11, 5).Value = '  -1.72%  '

$ws.Cells.Item(12, 4).Value = '1.817.04'
$ws.Cells.Item(12, 5).Value = '  -1.55%  '

$ws.Cells.Item(13, 4).NumberFormat = '@'
$ws.Cells.Item(13, 4).Value = '4.664'
$ws.Cells.Item(13, 5).Value = '  -3.16%  '

$ws.Cells.Item(14, 4).NumberFormat = '@'
$ws.Cells.Item(14, 4).Value = '0.6235'
$ws.Cells.Item(14, 5).Value = '  -2.23%  '

$ws.Cells.Item(15, 4).NumberFormat = '@'
$ws.Cells.Item(15, 4).Value = '0.000009297'
$ws.Cells.Item(15, 5).Value = '  -6.37%  '

$ws.Cells.Item(16, 4).NumberFormat = '@'
$ws.Cells.Item(16, 4).Value = '74.31'
$ws.Cells.Item(16, 5).Value = '  -6.80%  '

$ws.Cells.Item(17, 4).Value = '28.754.34'
$ws.Cells.Item(17, 5).Value = '  -1.63%  '

$ws.Cells.Item(18, 4).NumberFormat = '@'
$ws.Cells.Item(18, 4).Value = '5.409'
$ws.Cells.Item(18, 5).Value = '  -9.55%  '

$ws.Cells.Item(19, 5).Value = '  +0.51%  '

$ws.Cells.Item(20, 4).NumberFormat = '@'
$ws.Cells.Item(20, 4).Value = '207.28'
$ws.Cells.Item(20, 5).Value = '  -9.98%  '

$ws.Cells.Item(21, 4).NumberFormat = '@'
$ws.Cells.Item(21, 4).Value = '11.35'
$ws.Cells.Item(21, 5).Value = '  -3.92%  '

$ws.Cells.Item(22, 4).NumberFormat = '@'
$ws.Cells.Item(22, 4).Value = '6.761'
$ws.Cells.Item(22, 5).Value = '  -3.91%  '

$ws.Cells.Item(23, 4).NumberFormat = '@'
$ws.Cells.Item(23, 4).Value = '1.005'
$ws.Cells.Item(23, 5).Value = '  +0.34%  '

$ws.Cells.Item(24, 4).NumberFormat = '@'
$ws.Cells.Item(24, 4).Value = '153.91'
$ws.Cells.Item(24, 5).Value = '  -1.02%  '

$ws.Cells.Item(25, 2).Value = 'Stellar'
$ws.Cells.Item(25, 3).Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Cells.Item(25, 4).NumberFormat = '@'
$ws.Cells.Item(25, 4).Value = '0.1264'
$ws.Cells.Item(25, 5).Value = '  -2.52%  '

$ws.Cells.Item(26, 2).Value = 'Cosmos'
$ws.Cells.Item(26, 3).Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
$ws.Cells.Item(26, 4).NumberFormat = '@'
$ws.Cells.Item(26, 4).Value = '7.754'
$ws.Cells.Item(26, 5).Value = '  -4.20%  '

$ws.Cells.Item(27, 4).NumberFormat = '@'
$ws.Cells.Item(27, 4).Value = '16.23'
$ws.Cells.Item(27, 5).Value = '  -3.18%  '

$ws.Cells.Item(28, 2).Value = 'Hedera'
$ws.Cells.Item(28, 3).Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Cells.Item(28, 4).NumberFormat = '@'
$ws.Cells.Item(28, 4).Value = '0.06361'
$ws.Cells.Item(28, 5).Value = '  -5.33%  '

$ws.Cells.Item(29, 2).Value = 'Toncoin'
$ws.Cells.Item(29, 3).Value = 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton'
$ws.Cells.Item(29, 4).NumberFormat = '@'
$ws.Cells.Item(29, 4).Value = '1.405'
$ws.Cells.Item(29, 5).Value = '  -5.36%  '

$ws.Cells.Item(30, 4).NumberFormat = '@'
$ws.Cells.Item(30, 4).Value = '1.429'
$ws.Cells.Item(30, 5).Value = '  -2.05%  '

$ws.Cells.Item(31, 4).NumberFormat = '@'
$ws.Cells.Item(31, 4).Value = '3.702'
$ws.Cells.Item(31, 5).Value = '  -3.08%  '

$ws.Cells.Item(32, 4).NumberFormat = '@'
$ws.Cells.Item(32, 4).Value = '3.665'
$ws.Cells.Item(32, 5).Value = '  -4.66%  '

$ws.Cells.Item(33, 4).NumberFormat = '@'
$ws.Cells.Item(33, 4).Value = '1.677'
$ws.Cells.Item(33, 5).Value = '  -3.19%  '

$ws.Cells.Item(34, 4).NumberFormat = '@'
$ws.Cells.Item(34, 4).Value = '1.045'
$ws.Cells.Item(34, 5).Value = '  -7.04%  '

$ws.Cells.Item(35, 4).NumberFormat = '@'
$ws.Cells.Item(35, 4).Value = '2.533'
$ws.Cells.Item(35, 5).Value = '  -0.85%  '

$ws.Cells.Item(36, 4).NumberFormat = '@'
$ws.Cells.Item(36, 4).Value = '0.6313'
$ws.Cells.Item(36, 5).Value = '  -4.16%  '

$ws.Cells.Item(37, 4).NumberFormat = '@'
$ws.Cells.Item(37, 4).Value = '2.750'
$ws.Cells.Item(37, 5).Value = '  -0.23%  '

$ws.Cells.Item(38, 4).NumberFormat = '@'
$ws.Cells.Item(38, 4).Value = '6.419'
$ws.Cells.Item(38, 5).Value = '  -2.91%  '

$ws.Cells.Item(39, 4).NumberFormat = '@'
$ws.Cells.Item(39, 4).Value = '0.01698'
$ws.Cells.Item(39, 5).Value = '  -3.87%  '

$ws.Cells.Item(40, 4).Value = '1.134.71'
$ws.Cells.Item(40, 5).Value = '  -7.99%  '

$ws.Cells.Item(41, 4).NumberFormat = '@'
$ws.Cells.Item(41, 4).Value = '0.8678'
$ws.Cells.Item(41, 5).Value = '  -7.16%  '

$ws.Cells.Item(42, 4).NumberFormat = '@'
$ws.Cells.Item(42, 4).Value = '1.005'
$ws.Cells.Item(42, 5).Value = '  +0.46%  '

$ws.Cells.Item(43, 4).Value = '1.967.38'
$ws.Cells.Item(43, 5).Value = '  -1.13%  '

$ws.Cells.Item(44, 4).NumberFormat = '@'
$ws.Cells.Item(44, 4).Value = '99.58'
$ws.Cells.Item(44, 5).Value = '  -1.27%  '

$ws.Cells.Item(45, 4).NumberFormat = '@'
$ws.Cells.Item(45, 4).Value = '59.95'
$ws.Cells.Item(45, 5).Value = '  -5.69%  '

$ws.Cells.Item(46, 5).Value = '  -4.01%  '

$ws.Cells.Item(47, 2).Value = 'RenderToken'
$ws.Cells.Item(47, 3).Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Cells.Item(47, 4).NumberFormat = '@'
$ws.Cells.Item(47, 4).Value = '1.571'
$ws.Cells.Item(47, 5).Value = '  -3.84%  '

$ws.Cells.Item(48, 4).NumberFormat = '@'
$ws.Cells.Item(48, 4).Value = '0.4527'
$ws.Cells.Item(48, 5).Value = '  -0.80%  '

$ws.Cells.Item(49, 2).Value = 'Cronos'
$ws.Cells.Item(49, 3).Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Cells.Item(49, 4).NumberFormat = '@'
$ws.Cells.Item(49, 4).Value = '0.05479'
$ws.Cells.Item(49, 5).Value = '  -1.63%  '

$ws.Cells.Item(50, 4).NumberFormat = '@'
$ws.Cells.Item(50, 4).Value = '1.014'
$ws.Cells.Item(50, 5).Value = '  +0.98%  '

$ws.Cells.Item(51, 5).Value = '  -4.41%  '
